$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Năm" (Year) column from "2022-2023" to "2021-2022" for data rows 2-5
$ws.Range("D2").Value = "2021-2022"
$ws.Range("D3").Value = "2021-2022"
$ws.Range("D4").Value = "2021-2022"
$ws.Range("D5").Value = "2021-2022"

# Update the active selection to A2
$ws.Range("A2").Select()
